$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51 (pushes existing rows 51..82 down to 52..83)
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new record
$ws.Cells.Item(51, 1).Value = 9
$ws.Cells.Item(51, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(51, 3).Value = "Metropolitana"
$ws.Cells.Item(51, 4).Value = 44873
$ws.Cells.Item(51, 5).Value = 13
$ws.Cells.Item(51, 6).Value = 100112029
$ws.Cells.Item(51, 7).Value = "Orégano"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 11
$ws.Cells.Item(51, 11).Value = 15000
$ws.Cells.Item(51, 12).Value = 18000
$ws.Cells.Item(51, 13).Value = 16636
$ws.Cells.Item(51, 14).Value = "`$/docena de atados"
$ws.Cells.Item(51, 15).Value = "Región Metropolitana"
$ws.Cells.Item(51, 16).Value = 5545
$ws.Cells.Item(51, 17).Value = 3
$ws.Cells.Item(51, 18).Value = "Hortaliza"
